$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (1-indexed columns A..Y) - target stored <col> width values.
# The engine adds a fixed 5/6 (0.8333333333333334) padding on top of whatever
# ColumnWidth we set, so compensate by subtracting it here.
$pad = 0.8333333333333334
$colWidths = @(26,15,21,19,17,9,14,22,21,13,7,17,15,5,5,9,9,13,15,21,13,13,14,14,18)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $col = $i + 1
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$i] - $pad
}

# New header cells I1:Y1 (values first)
$headers = @{
    "I1" = "Datetime"
    "J1" = "Ip Address"
    "K1" = "Model"
    "L1" = "Hospital Center"
    "M1" = "Protocol Code"
    "N1" = "Age"
    "O1" = "Sex"
    "P1" = "Max Dim"
    "Q1" = "Min Dim"
    "R1" = "Veinous Inf"
    "S1" = "Arterious Inf"
    "T1" = "Duct Ret/Ductal Inv"
    "U1" = "Vessel Comp"
    "V1" = "Reg Margins"
    "W1" = "Echogenicity"
    "X1" = "Mult Lesions"
    "Y1" = "Prediction"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (A1) onto the newly added header cells I1:Y1.
$ws.Range("A1").Copy()
$ws.Range("I1:Y1").PasteSpecial(-4122)

# Row 19
$ws.Range("E19").Value = 0
$ws.Range("I19").Value = "2025-05-13 16:28:54"
$ws.Range("J19").Value = "93.66.99.95"
$ws.Range("K19").Value = "DT"
$ws.Range("O19").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("T19").Value = 0
$ws.Range("U19").Value = 0
$ws.Range("V19").Value = 0
$ws.Range("W19").Value = 0
$ws.Range("X19").Value = 0
$ws.Range("Y19").Value = "22.95% Malignant"

# Row 20
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = "2025-05-13 16:29:49"
$ws.Range("J20").Value = "93.66.99.95"
$ws.Range("K20").Value = "DT"
$ws.Range("O20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("U20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("W20").Value = 0
$ws.Range("X20").Value = 0
$ws.Range("Y20").Value = "22.95% Malignant"
